$d = $word.ActiveDocument
$d.Content.Find.Execute("table name, 3, 4 FROM information_schema.tables WHERE table schema", $true, $false, $false, $false, $false, $true, 1, $false, "table_name, 3, 4 FROM information_schema.tables WHERE table_schema", 2)
